# Add data for 2022-07-06
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the "as of" date references
$ws.Name = "Through 2022-06-28"
$ws.Range("I1").Value = "2022 (through 06-28)"

# Update the June (row 7) and Total (row 14) figures for the 2022 column (I)
$ws.Range("I7").Value = 133
$ws.Range("I14").Value = 796
